$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 (pushes old rows 4-7 down to 5-8),
# so PREPROD rows become rows 5-8 and QA rows occupy 2-4.
$ws.Rows.Item(4).Insert()

# New row's A4/B4 should carry QA formatting/content
$ws.Cells.Item(4, 1).Value = "QA"

$ws.Cells.Item(2, 2).Value = "'1120194100403 "
$ws.Cells.Item(3, 2).Value = "'1220194200660"
$ws.Cells.Item(4, 2).Value = "'0420194406623"

# Set selection to B5 as in diff
$ws.Range("B5").Select()
